$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "timestamp" column (O) for all data rows (2 through 398)
# from the old value "2023-01-07 12:54:35" to the new value "2023-01-07 20:49:26"
$ws.Range("O2:O398").Value = "2023-01-07 20:49:26"
